# Update cryptos list values (prices and 1h volume percentages)
# Applies updated market data snapshot to the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.145.97"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.462.28"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.25"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.90"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.461.03"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.45"
$ws.Range("E14").Value = "  +7.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +5.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.898.77"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.122.97"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.461.44"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.97"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.09"
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.99"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.14"
$ws.Range("E23").Value = "  +10.72%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.26"
$ws.Range("E25").Value = "  +25.66%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.44"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "652.13"
$ws.Range("E27").Value = "  +9.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("E29").Value = "  +7.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.585.62"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.22"
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  +5.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("E34").Value = "  +4.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₆0393"
$ws.Range("E35").Value = "  +39.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +4.07%  "
$ws.Range("E39").Value = "  +5.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.375"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.88"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "152.38"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("E43").Value = "  +9.78%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.00"
$ws.Range("E47").Value = "  +27.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.25"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.64"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.76"
$ws.Range("E50").Value = "  +4.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.609"
$ws.Range("E51").Value = "  +2.71%  "
